$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "ruhi"
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "31413"
$ws.Cells.Item(4, 4).Value = "sad"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "2"
$ws.Cells.Item(4, 7).Value = "ÖDEME ALINDI"
$ws.Cells.Item(4, 8).Value = "26-03-2023"
